$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SanityTC")

# --- New header cells CD1:CF1 (string values, same highlighted style as CC1) ---
$ws.Range("CC1").Copy()
$ws.Range("CD1:CF1").PasteSpecial(-4122)

$ws.Range("CD1").Value = "demographics"
$ws.Range("CE1").Value = "subquestions"
$ws.Range("CF1").Value = "symbolratingscale"

# --- New blank data cells CD2:CF3 (bordered, same look as column CC) ---
foreach ($addr in @("CD2","CE2","CF2","CD3","CE3","CF3")) {
    $c = $ws.Range($addr)
    $c.Borders.Item(7).LineStyle = 1
    $c.Borders.Item(7).Weight = 2
    $c.Borders.Item(7).ColorIndex = 1
    $c.Borders.Item(8).LineStyle = 1
    $c.Borders.Item(8).Weight = 2
    $c.Borders.Item(8).ColorIndex = 1
    $c.Borders.Item(9).LineStyle = 1
    $c.Borders.Item(9).Weight = 2
    $c.Borders.Item(9).ColorIndex = 1
    $c.Borders.Item(10).LineStyle = 1
    $c.Borders.Item(10).Weight = 2
    $c.Borders.Item(10).ColorIndex = 1
}

# --- Move the active selection to CD8, matching the saved view state ---
[void]$ws.Range("CD8").Select()
